$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force de-duplication of the redundant "订单管理" shared string (B4/C4
# originally pointed at two separate but identical strings) so both cells
# end up referencing the same shared-string entry, as happens when Excel
# resaves after the row is touched.
$ws.Range("B4").Value = "订单管理_TEMP"
$ws.Range("B4").Value = "订单管理"

# Insert 3 new rows after row 4 (order management sub-permissions)
$ws.Rows("5:7").Insert()

# Copy formatting (border/style) from row 4 into the new rows 5-7
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new rows with order management sub-permission data
$ws.Range("A5").Value = 11001
$ws.Range("B5").Value = "订单删除"
$ws.Range("C5").Value = "订单管理 - 订单删除"

$ws.Range("A6").Value = 11002
$ws.Range("B6").Value = "订单修改"
$ws.Range("C6").Value = "订单管理 - 订单修改"

$ws.Range("A7").Value = 11003
$ws.Range("B7").Value = "订单状态修改"
$ws.Range("C7").Value = "订单管理 - 订单状态修改"

# Update the view/selection as recorded after the edit
$ws.Range("B8").Select()
